$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nppc"
$ws.Range("C2").Value = "Npr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01599966666666667
$ws.Range("H2").Value = 0.047999
$ws.Range("I2").Value = 0.002186576436924975
$ws.Range("J2").Value = 0.002186576436924975
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.243779
$ws.Range("N2").Value = 48.731337
$ws.Range("O2").Value = 0.4875086269552682
$ws.Range("P2").Value = 0.4875086269552682
$ws.Range("Q2").Value = 0.259895049407
$ws.Range("R2").Value = 2.339055444663
$ws.Range("S2").Value = 0.001065974876498037
$ws.Range("T2").Value = 0.001065974876498037

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nppc"
$ws.Range("C3").Value = "Npr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01599966666666667
$ws.Range("H3").Value = 0.047999
$ws.Range("I3").Value = 0.002186576436924975
$ws.Range("J3").Value = 0.002186576436924975
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 13.462409
$ws.Range("N3").Value = 40.387227
$ws.Range("O3").Value = 0.4040340937352229
$ws.Range("P3").Value = 0.4040340937352229
$ws.Range("Q3").Value = 0.2153940565303333
$ws.Range("R3").Value = 1.938546508773
$ws.Range("S3").Value = 0.0008834514290757751
$ws.Range("T3").Value = 0.0008834514290757751

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nppc"
$ws.Range("C4").Value = "Npr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01599966666666667
$ws.Range("H4").Value = 0.047999
$ws.Range("I4").Value = 0.002186576436924975
$ws.Range("J4").Value = 0.002186576436924975
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02644566666666667
$ws.Range("N4").Value = 0.079337
$ws.Range("O4").Value = 0.000793687887872851
$ws.Range("P4").Value = 0.000793687887872851
$ws.Range("Q4").Value = 0.0004231218514444445
$ws.Range("R4").Value = 0.003808096663
$ws.Range("S4").Value = 0.000001735459233895528
$ws.Range("T4").Value = 0.000001735459233895528

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Nppc"
$ws.Range("C5").Value = "Npr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01599966666666667
$ws.Range("H5").Value = 0.047999
$ws.Range("I5").Value = 0.002186576436924975
$ws.Range("J5").Value = 0.002186576436924975
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.587349
$ws.Range("N5").Value = 10.762047
$ws.Range("O5").Value = 0.1076635914216362
$ws.Range("P5").Value = 0.1076635914216362
$ws.Range("Q5").Value = 0.05739638821699999
$ws.Range("R5").Value = 0.5165674939529999
$ws.Range("S5").Value = 0.0002354146721172676
$ws.Range("T5").Value = 0.0002354146721172676

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nppc"
$ws.Range("C6").Value = "Npr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.899989666666667
$ws.Range("H6").Value = 5.699969
$ws.Range("I6").Value = 0.2596599493031691
$ws.Range("J6").Value = 0.2596599493031691
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.243779
$ws.Range("N6").Value = 48.731337
$ws.Range("O6").Value = 0.4875086269552682
$ws.Range("P6").Value = 0.4875086269552682
$ws.Range("Q6").Value = 30.863012247617
$ws.Range("R6").Value = 277.767110228553
$ws.Range("S6").Value = 0.1265864653600625
$ws.Range("T6").Value = 0.1265864653600625

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nppc"
$ws.Range("C7").Value = "Npr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.899989666666667
$ws.Range("H7").Value = 5.699969
$ws.Range("I7").Value = 0.2596599493031691
$ws.Range("J7").Value = 0.2596599493031691
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.462409
$ws.Range("N7").Value = 40.387227
$ws.Range("O7").Value = 0.4040340937352229
$ws.Range("P7").Value = 0.4040340937352229
$ws.Range("Q7").Value = 25.57843798844033
$ws.Range("R7").Value = 230.205941895963
$ws.Range("S7").Value = 0.1049114722960398
$ws.Range("T7").Value = 0.1049114722960398

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Nppc"
$ws.Range("C8").Value = "Npr2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.899989666666667
$ws.Range("H8").Value = 5.699969
$ws.Range("I8").Value = 0.2596599493031691
$ws.Range("J8").Value = 0.2596599493031691
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02644566666666667
$ws.Range("N8").Value = 0.079337
$ws.Range("O8").Value = 0.000793687887872851
$ws.Range("P8").Value = 0.000793687887872851
$ws.Range("Q8").Value = 0.05024649339477778
$ws.Range("R8").Value = 0.452218440553
$ws.Range("S8").Value = 0.0002060889567276039
$ws.Range("T8").Value = 0.0002060889567276039

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Nppc"
$ws.Range("C9").Value = "Npr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.899989666666667
$ws.Range("H9").Value = 5.699969
$ws.Range("I9").Value = 0.2596599493031691
$ws.Range("J9").Value = 0.2596599493031691
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.587349
$ws.Range("N9").Value = 10.762047
$ws.Range("O9").Value = 0.1076635914216362
$ws.Range("P9").Value = 0.1076635914216362
$ws.Range("Q9").Value = 6.815926030727
$ws.Range("R9").Value = 61.343334276543
$ws.Range("S9").Value = 0.02795592269033917
$ws.Range("T9").Value = 0.02795592269033917

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Nppc"
$ws.Range("C10").Value = "Npr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.401233333333334
$ws.Range("H10").Value = 16.2037
$ws.Range("I10").Value = 0.738153474259906
$ws.Range("J10").Value = 0.738153474259906
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 16.243779
$ws.Range("N10").Value = 48.731337
$ws.Range("O10").Value = 0.4875086269552682
$ws.Range("P10").Value = 0.4875086269552682
$ws.Range("Q10").Value = 87.7364405941
$ws.Range("R10").Value = 789.6279653469001
$ws.Range("S10").Value = 0.3598561867187077
$ws.Range("T10").Value = 0.3598561867187077

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Nppc"
$ws.Range("C11").Value = "Npr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5.401233333333334
$ws.Range("H11").Value = 16.2037
$ws.Range("I11").Value = 0.738153474259906
$ws.Range("J11").Value = 0.738153474259906
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.462409
$ws.Range("N11").Value = 40.387227
$ws.Range("O11").Value = 0.4040340937352229
$ws.Range("P11").Value = 0.4040340937352229
$ws.Range("Q11").Value = 72.71361223776667
$ws.Range("R11").Value = 654.4225101398999
$ws.Range("S11").Value = 0.2982391700101072
$ws.Range("T11").Value = 0.2982391700101072

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Nppc"
$ws.Range("C12").Value = "Npr2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 5.401233333333334
$ws.Range("H12").Value = 16.2037
$ws.Range("I12").Value = 0.738153474259906
$ws.Range("J12").Value = 0.738153474259906
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02644566666666667
$ws.Range("N12").Value = 0.079337
$ws.Range("O12").Value = 0.000793687887872851
$ws.Range("P12").Value = 0.000793687887872851
$ws.Range("Q12").Value = 0.1428392163222222
$ws.Range("R12").Value = 1.2855529469
$ws.Range("S12").Value = 0.0005858634719113516
$ws.Range("T12").Value = 0.0005858634719113516

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Nppc"
$ws.Range("C13").Value = "Npr2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 5.401233333333334
$ws.Range("H13").Value = 16.2037
$ws.Range("I13").Value = 0.738153474259906
$ws.Range("J13").Value = 0.738153474259906
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.587349
$ws.Range("N13").Value = 10.762047
$ws.Range("O13").Value = 0.1076635914216362
$ws.Range("P13").Value = 0.1076635914216362
$ws.Range("Q13").Value = 19.3761089971
$ws.Range("R13").Value = 174.3849809739
$ws.Range("S13").Value = 0.07947225405917976
$ws.Range("T13").Value = 0.07947225405917976
